# Update status task in sprint
# Row 5 ("MODELLING WITH NO IMAGE AUGMENTATION") is updated:
#   - TANGGAL SELESAI (column D) is filled in with the completion date
#   - STATUS (column F) is changed from WAITING to DONE

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D5").Value = "8 Nopember 2021"
$ws.Range("F5").Value = "DONE"

# Update the active cell selection on the sheet to reflect where the user
# left off editing.
$ws.Activate()
$ws.Range("F6").Select()
